# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values for each observation row (rows 2-11).
# Recalculated K values replace the previous Strike#-derived values.
$kValues = @{
    2  = 3
    3  = 2
    4  = 1
    5  = 2
    6  = 1
    7  = 3
    8  = 4
    9  = 0
    10 = 3
    11 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
